$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "hockey knee pads"
$ws.Range("A2").Value = "knee black leggings"
$ws.Range("A3").Value = "knee compression leggings"
$ws.Range("A4").Value = "knee compression pad"
$ws.Range("A5").Value = "knee compression running"
$ws.Range("A6").Value = "knee gel pad"
$ws.Range("A7").Value = "knee pad for men"
$ws.Range("A8").Value = "knee pad for volleyball"
$ws.Range("A9").Value = "knee pad for workout"
$ws.Range("A10").Value = "knee pad leggings"
$ws.Range("A11").Value = "knee pads athletic"
$ws.Range("A12").Value = "knee pads basketball youth"
$ws.Range("A13").Value = "knee pads black"
$ws.Range("A14").Value = "knee pads boys"
$ws.Range("A15").Value = "knee pads extra large"
$ws.Range("A16").Value = "knee pads for football"
$ws.Range("A17").Value = "knee pads for gym men"
$ws.Range("A18").Value = "knee pads for men"
$ws.Range("A19").Value = "knee pads for running men"
$ws.Range("A20").Value = "knee pads for soccer"
$ws.Range("A21").Value = "knee pads honeycomb"
$ws.Range("A22").Value = "knee pads mens"
$ws.Range("A23").Value = "knee pads skating youth"
$ws.Range("A24").Value = "knee pads snowboarding"
$ws.Range("A25").Value = "knee pads soccer"
$ws.Range("A26").Value = "knee pads sport"
$ws.Range("A27").Value = "knee pads squat"
$ws.Range("A28").Value = "knee pads thick"
$ws.Range("A29").Value = "knee pads weightlifting"
$ws.Range("A30").Value = "knee pads youth"
$ws.Range("A31").Value = "knee pants for men"
$ws.Range("A32").Value = "knee protection"
$ws.Range("A33").Value = "knee protector gym"
$ws.Range("A34").Value = "knee protectors"
$ws.Range("A35").Value = "knee replacement aids"
$ws.Range("A36").Value = "knee support gym"
$ws.Range("A37").Value = "knee support pants"
$ws.Range("A38").Value = "knee support pants men"
$ws.Range("A39").Value = "kneepads for volleyball"
$ws.Range("A40").Value = "leg compression pants"
$ws.Range("A41").Value = "leg compression tights"
$ws.Range("A42").Value = "leg protectors for men"
$ws.Range("A43").Value = "legging for basketball boys"
$ws.Range("A44").Value = "leggings for men sport gym"
$ws.Range("A45").Value = "leggings knee"
$ws.Range("A46").Value = "leggings medium"
$ws.Range("A47").Value = "leggings mens"
$ws.Range("A48").Value = "leggins training"
$ws.Range("A49").Value = "lightweight athletic pants for men"
$ws.Range("A50").Value = "lightweight pants men"
$ws.Range("A51").Value = "lightweight sports pants men"
$ws.Range("A52").Value = "mcdavid basketball knee pads 6446"
$ws.Range("A53").Value = "mcdavid basketball knee pads black"
$ws.Range("A54").Value = "men basketball tights"
$ws.Range("A55").Value = "men capri pants"
$ws.Range("A56").Value = "men compression pants"
$ws.Range("A57").Value = "men leggings pack"
$ws.Range("A58").Value = "men leggings running"
$ws.Range("A59").Value = "men leggings tall"
$ws.Range("A60").Value = "men running tights nike"
$ws.Range("A61").Value = "men tights legging"
$ws.Range("A62").Value = "mens athletic compression pants"
$ws.Range("A63").Value = "mens athletic leggings black"
$ws.Range("A64").Value = "mens basketball knee pads"
$ws.Range("A65").Value = "mens compression 3 4 pants"
$ws.Range("A66").Value = "mens compression leggings 3 4"
$ws.Range("A67").Value = "mens compression tights pants"
$ws.Range("A68").Value = "mens cycling tights"
$ws.Range("A69").Value = "mens hiking pants lightweight"
$ws.Range("A70").Value = "mens jogging tights"
$ws.Range("A71").Value = "mens lacrosse pads"
$ws.Range("A72").Value = "mens leggings"
$ws.Range("A73").Value = "mens lightweight workout pants"
$ws.Range("A74").Value = "mens running tights green"
$ws.Range("A75").Value = "mens running tights orange"
$ws.Range("A76").Value = "mens running tights yellow"
$ws.Range("A77").Value = "mens soccer clothing"
$ws.Range("A78").Value = "mens stretch thermal pants"
$ws.Range("A79").Value = "mens swim leggings"
$ws.Range("A80").Value = "mens tights basketball"
$ws.Range("A81").Value = "mens tights capri"
$ws.Range("A82").Value = "mens tights leggings"
$ws.Range("A83").Value = "mens underarmour snow pants"
$ws.Range("A84").Value = "mens workout pants"
$ws.Range("A85").Value = "multicam pants with knee pads"
$ws.Range("A86").Value = "nike basketball pads"
$ws.Range("A87").Value = "nike volleyball knee pads youth girls"
$ws.Range("A88").Value = "padded compression"
$ws.Range("A89").Value = "padded soccer pants"
$ws.Range("A90").Value = "pain in thigh joint"
$ws.Range("A91").Value = "pant knee pad inserts"
$ws.Range("A92").Value = "pantalon con rodilleras"
$ws.Range("A93").Value = "pants with knee pads"
$ws.Range("A94").Value = "polyester pants men pants"
$ws.Range("A95").Value = "protective knee pads for men"
$ws.Range("A96").Value = "rash guard men bjj"
$ws.Range("A97").Value = "reebok compression pants"
$ws.Range("A98").Value = "running compression leg"
$ws.Range("A99").Value = "running knee"
$ws.Range("A100").Value = "running pants youth"
